$d = $word.ActiveDocument

# --- 1. Restart the numbering for the Problem #2 "A." list item (numId 8 -> new list, numId 9) ---
$pA = $d.Paragraphs.Item(23)
$sourceTemplate = $d.ListTemplates.Item(1)
$pA.Range.ListFormat.ApplyListTemplateWithLevel($sourceTemplate, $false, 2, $false, 1)

# Bring the freshly-minted list's sub-levels back in line with the rest of the
# document's numbered lists (decimal / lowerLetter / lowerRoman repeating pattern).
$newTemplate = $pA.Range.ListFormat.ListTemplate
$newTemplate.ListLevels.Item(2).NumberStyle = 4   # lowerLetter
$newTemplate.ListLevels.Item(3).NumberStyle = 2   # lowerRoman
$newTemplate.ListLevels.Item(5).NumberStyle = 4   # lowerLetter
$newTemplate.ListLevels.Item(6).NumberStyle = 2   # lowerRoman
$newTemplate.ListLevels.Item(8).NumberStyle = 4   # lowerLetter
$newTemplate.ListLevels.Item(9).NumberStyle = 2   # lowerRoman

# --- 2. Add the new "sub goals" paragraphs for problem #2, just above the
#        trailing bookmark paragraph. ---
$n = $d.Paragraphs.Count
$bookmarkPara = $d.Paragraphs.Item($n)

$bookmarkPara.Range.InsertParagraphBefore()
$bookmarkPara.Range.InsertParagraphBefore()

$p1 = $d.Paragraphs.Item($n)
$p1.Style = "Normal"
$p1.Range.Text = "2. A. the sub goals are to count each pair"

$p2 = $d.Paragraphs.Item($n + 1)
$p2.Style = "Normal"
$p2.Range.Text = "          Try to use actual socks as a demonstration."

# --- 3. The old bookmark paragraph loses its ListParagraph style and gets a
#        leading space run in front of the (still present) bookmark. ---
$bookmarkPara2 = $d.Paragraphs.Item($n + 2)
$bookmarkPara2.Style = "Normal"
$bookmarkPara2.Range.InsertBefore(" ")

# --- 4. A final, trailing empty paragraph is appended at the very end. ---
$d.Content.InsertParagraphAfter()
$trailing = $d.Paragraphs.Item($d.Paragraphs.Count)
$trailing.Style = "Normal"
